$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6
$ws.Range("J6").Value = 1.05
$ws.Range("L6").Value = 1.33

# Row 11
$ws.Range("G11").Value = 1.85
$ws.Range("H11").Value = 3.4
$ws.Range("I11").Value = 3.9
$ws.Range("L11").Value = 1.3
$ws.Range("M11").Value = 2.95
$ws.Range("N11").Value = 1.88
$ws.Range("O11").Value = 1.72
$ws.Range("P11").Value = 1.39
$ws.Range("Q11").Value = 2.57
$ws.Range("R11").Value = 1.78
$ws.Range("S11").Value = 1.83
$ws.Range("T11").Value = 7
$ws.Range("U11").Value = 8.75
$ws.Range("V11").Value = 8.25
$ws.Range("W11").Value = 15.5
$ws.Range("X11").Value = 15
$ws.Range("Y11").Value = 27
$ws.Range("Z11").Value = 9.5
$ws.Range("AA11").Value = 6.6
$ws.Range("AB11").Value = 15.5
$ws.Range("AC11").Value = 75
$ws.Range("AD11").Value = 600
$ws.Range("AE11").Value = 10.5
$ws.Range("AF11").Value = 21
$ws.Range("AI11").Value = 37
$ws.Range("AJ11").Value = 45

# Row 13
$ws.Range("G13").Value = 2.4
$ws.Range("H13").Value = 3.4
$ws.Range("I13").Value = 2.7
$ws.Range("N13").Value = 2.05
$ws.Range("O13").Value = 1.75
$ws.Range("P13").Value = 1.4
$ws.Range("Q13").Value = 2.75
$ws.Range("T13").Value = 8
$ws.Range("Y13").Value = 29
$ws.Range("Z13").Value = 9.5
$ws.Range("AD13").Value = 251
$ws.Range("AE13").Value = 8.5
